$d = $word.ActiveDocument

# The header title "Techniken der Empirischen Parlamentsanalyse in R" was
# split across two runs ("T" + "echniken ... in R") for no semantic reason.
# Re-typing it as a single Find/Replace over the header range merges it
# back into one run, cleaning up the file structure without changing the
# visible text or formatting.
$header = $d.Sections(1).Headers(1)
$header.Range.Find.Execute("Techniken der Empirischen Parlamentsanalyse in R", $true, $false, $false, $false, $false, $true, 1, $false, "Techniken der Empirischen Parlamentsanalyse in R", 2)
